# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between the existing "2021-Q4" sheet
#    and the "总计" (totals) sheet, and fill it with the per-fund holding
#    data for 2022-Q1.
# 2. Update the "总计" sheet: insert a new row for "2022-Q1" above the
#    existing "2021-Q4" row (newest quarter on top).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet right before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Re-fetch sheets by name from now on - inserting/renaming sheets shifts
# tab positions, and an Item(<index>) reference grabbed beforehand would
# silently start pointing at whatever sheet now sits at that index.
$q4 = $wb.Worksheets.Item("2021-Q4")   # formatting template
$q1 = $wb.Worksheets.Item("2022-Q1")

# Copy the header-row formatting (bold font + border) from "2021-Q4"
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q4.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B (fund code) and D:G (numeric-looking text) must stay text,
# same as how the data was originally exported.
$q1.Range("B2:B9").NumberFormat = "@"
$q1.Range("D2:G9").NumberFormat = "@"

$rows = @(
    @{ idx = 0; code = "700001"; name = "平安行业先锋混合";             scale = "2.32"; pos = "91.67"; ratio = "2.99"; value = "0.0694"; rank = 9 },
    @{ idx = 1; code = "001664"; name = "平安鑫安混合A";               scale = "0.86"; pos = "29.46"; ratio = "1.04"; value = "0.0089"; rank = 7 },
    @{ idx = 2; code = "011761"; name = "平安鑫瑞混合型证券投资基金A"; scale = "1.09"; pos = "20.41"; ratio = "0.52"; value = "0.0057"; rank = 9 },
    @{ idx = 3; code = "007049"; name = "平安鑫安混合E";               scale = "0.50"; pos = "29.46"; ratio = "1.04"; value = "0.0052"; rank = 7 },
    @{ idx = 4; code = "003626"; name = "平安鑫利灵活配置混合A";       scale = "0.41"; pos = "27.68"; ratio = "0.87"; value = "0.0036"; rank = 9 },
    @{ idx = 5; code = "006433"; name = "平安鑫利灵活配置混合C";       scale = "0.25"; pos = "27.68"; ratio = "0.87"; value = "0.0022"; rank = 9 },
    @{ idx = 6; code = "011762"; name = "平安鑫瑞混合型证券投资基金C"; scale = "0.34"; pos = "20.41"; ratio = "0.52"; value = "0.0018"; rank = 9 },
    @{ idx = 7; code = "001665"; name = "平安鑫安混合C";               scale = "0.02"; pos = "29.46"; ratio = "1.04"; value = "0.0002"; rank = 7 }
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row.idx
    $q1.Cells.Item($r, 2).Value = $row.code
    $q1.Cells.Item($r, 3).Value = $row.name
    $q1.Cells.Item($r, 4).Value = $row.scale
    $q1.Cells.Item($r, 5).Value = $row.pos
    $q1.Cells.Item($r, 6).Value = $row.ratio
    $q1.Cells.Item($r, 7).Value = $row.value
    $q1.Cells.Item($r, 8).Value = $row.rank
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: add the 2022-Q1 summary row above 2021-Q4,
#    pushing the existing 2021-Q4 row down from row 2 to row 3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Inserting a row copies formatting down from the row above (the bold
# header row) - clear it so the new data row matches the plain,
# unstyled look of the other "总计" data rows.
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.1

# The old row 2 ("2021-Q4") shifted down to row 3; re-number its index cell.
$total.Cells.Item(3, 1).Value = 1

# Copy the "A" column summary-sheet style (bold/bordered) onto the newly
# inserted index cell, matching the existing row's formatting.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Restore the originally active tab/selection ("2021-Q4").
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Activate()
